$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("I2").Value = 2.88
$ws.Range("L2").Value = 3.75
$ws.Range("W2").Value = 6
$ws.Range("AH2").Value = 11
$ws.Range("AJ2").Value = 29
$ws.Range("AQ2").Value = 81
$ws.Range("G6").Value = 1.44
$ws.Range("G7").Value = 1.91
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 4.5
$ws.Range("J7").Value = 2.63
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 4.75
$ws.Range("W7").Value = 6
$ws.Range("X7").Value = 8
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 15
$ws.Range("AA7").Value = 17
$ws.Range("AC7").Value = 7.5
$ws.Range("AF7").Value = 67
$ws.Range("AG7").Value = 10
$ws.Range("AH7").Value = 21
$ws.Range("AI7").Value = 15
$ws.Range("AK7").Value = 41
$ws.Range("AM7").Value = 1250
$ws.Range("AN7").Value = 3.75
$ws.Range("AO7").Value = 11
$ws.Range("AP7").Value = 23
$ws.Range("AW7").Value = 6
$ws.Range("AX7").Value = 23
$ws.Range("AZ7").Value = 81
$ws.Range("BA7").Value = 126
$ws.Range("G8").Value = 1.67
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 5.75
$ws.Range("J8").Value = 2.38
$ws.Range("K8").Value = 2.05
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("X8").Value = 7
$ws.Range("Z8").Value = 12
$ws.Range("AD8").Value = 6.5
$ws.Range("AE8").Value = 19
$ws.Range("AG8").Value = 12
$ws.Range("AH8").Value = 26
$ws.Range("AI8").Value = 19
$ws.Range("AJ8").Value = 51
$ws.Range("AN8").Value = 3.5
$ws.Range("AO8").Value = 9
$ws.Range("AX8").Value = 29
$ws.Range("Q11").Value = 2.6
$ws.Range("R11").Value = 1.48
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 10
$ws.Range("M14").Value = 1.1
$ws.Range("N14").Value = 7
$ws.Range("O14").Value = 1.5
$ws.Range("R14").Value = 1.5
$ws.Range("G18").Value = 2.75
$ws.Range("I18").Value = 2.42
$ws.Range("J18").Value = 3.35
$ws.Range("L18").Value = 3.05
$ws.Range("N18").Value = 5.95
$ws.Range("W18").Value = 7.2
$ws.Range("X18").Value = 12.5
$ws.Range("Y18").Value = 10.75
$ws.Range("Z18").Value = 32
$ws.Range("AA18").Value = 27
$ws.Range("AC18").Value = 7.5
$ws.Range("AD18").Value = 6.2
$ws.Range("AG18").Value = 6.6
$ws.Range("AH18").Value = 10.75
$ws.Range("AI18").Value = 10
$ws.Range("AJ18").Value = 25
$ws.Range("AK18").Value = 23
$ws.Range("AN18").Value = 4.5
$ws.Range("AO18").Value = 15
$ws.Range("AP18").Value = 26
$ws.Range("AQ18").Value = 70
$ws.Range("AS18").Value = 400
$ws.Range("AW18").Value = 4.15
$ws.Range("AX18").Value = 13
$ws.Range("AZ18").Value = 55
